$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1526.875
$ws.Range("I4").Value = 1526.875
$ws.Range("K4").Value = 1526.875
$ws.Range("M4").Value = -1412.875
$ws.Range("H86").Value = 5555.4443
$ws.Range("I86").Value = 4999.6665
$ws.Range("J86").Value = 5833.3335
$ws.Range("K86").Value = 4999.6665
$ws.Range("L86").Value = 5833.3335
$ws.Range("M86").Value = -3876.6665
$ws.Range("N86").Value = -8079.3335
$ws.Range("H88").Value = 6950.7
$ws.Range("J88").Value = 7056.3335
$ws.Range("L88").Value = 7056.3335
$ws.Range("N88").Value = -7868.3335
$ws.Range("H89").Value = 5555.4443
$ws.Range("I89").Value = 4999.6665
$ws.Range("J89").Value = 5833.3335
$ws.Range("K89").Value = 24998.3325
$ws.Range("L89").Value = 29166.6675
$ws.Range("M89").Value = -19382.3325
$ws.Range("N89").Value = -40398.6675
$ws.Range("H91").Value = 6950.7
$ws.Range("J91").Value = 7056.3335
$ws.Range("L91").Value = 7056.3335
$ws.Range("N91").Value = -9864.333500000001
$ws.Range("H106").Value = 29960.875
$ws.Range("I106").Value = 18445.666
$ws.Range("K106").Value = 18445.666
$ws.Range("M106").Value = -17814.666
$ws.Range("H107").Value = 3158.2727
$ws.Range("J107").Value = 5335.3335
$ws.Range("L107").Value = 5335.3335
$ws.Range("N107").Value = -9175.333500000001
$ws.Range("H116").Value = 7376.143
$ws.Range("I116").Value = 3664
$ws.Range("J116").Value = 8861
$ws.Range("K116").Value = 3664
$ws.Range("L116").Value = 8861
$ws.Range("M116").Value = -222
$ws.Range("N116").Value = -15745
$ws.Range("H132").Value = 2195.4348
$ws.Range("I132").Value = 1397.3636
$ws.Range("K132").Value = 4192.0908
$ws.Range("M132").Value = -1662.0908
$ws.Range("H138").Value = 2943.6135
$ws.Range("J138").Value = 3367.5217
$ws.Range("L138").Value = 10102.5651
$ws.Range("N138").Value = -20382.5651

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 15528
$ws.Range("I36").Value = 11950.4
$ws.Range("K36").Value = 11950.4
$ws.Range("M36").Value = -11604.4
$ws.Range("H132").Value = 4683.269
$ws.Range("I132").Value = 3436.475
$ws.Range("J132").Value = 8839.25
$ws.Range("K132").Value = 10309.425
$ws.Range("L132").Value = 26517.75
$ws.Range("M132").Value = -7779.424999999999
$ws.Range("N132").Value = -31577.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3542.2856
$ws.Range("I99").Value = 3435.7273
$ws.Range("J99").Value = 3933
$ws.Range("K99").Value = 3435.7273
$ws.Range("L99").Value = 3933
$ws.Range("M99").Value = -1937.7273
$ws.Range("N99").Value = -6929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30014
$ws.Range("I31").Value = 5739.8
$ws.Range("J31").Value = 33294.297
$ws.Range("K31").Value = 5739.8
$ws.Range("L31").Value = 33294.297
$ws.Range("M31").Value = -5444.8
$ws.Range("N31").Value = -33884.297
$ws.Range("H34").Value = 30014
$ws.Range("I34").Value = 5739.8
$ws.Range("J34").Value = 33294.297
$ws.Range("K34").Value = 5739.8
$ws.Range("L34").Value = 33294.297
$ws.Range("M34").Value = -5537.8
$ws.Range("N34").Value = -33698.297
$ws.Range("H62").Value = 9873.111000000001
$ws.Range("I62").Value = 7288.3335
$ws.Range("J62").Value = 11165.5
$ws.Range("K62").Value = 7288.3335
$ws.Range("L62").Value = 11165.5
$ws.Range("M62").Value = -6664.3335
$ws.Range("N62").Value = -12413.5
$ws.Range("H65").Value = 9873.111000000001
$ws.Range("I65").Value = 7288.3335
$ws.Range("J65").Value = 11165.5
$ws.Range("K65").Value = 36441.6675
$ws.Range("L65").Value = 55827.5
$ws.Range("M65").Value = -33321.6675
$ws.Range("N65").Value = -62067.5
$ws.Range("H86").Value = 6861
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("H89").Value = 6861
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384
$ws.Range("H99").Value = 3873.625
$ws.Range("I99").Value = 3747.75
$ws.Range("J99").Value = 3999.5
$ws.Range("K99").Value = 3747.75
$ws.Range("L99").Value = 3999.5
$ws.Range("M99").Value = -2249.75
$ws.Range("N99").Value = -6995.5
$ws.Range("H126").Value = 3873.625
$ws.Range("I126").Value = 3747.75
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 11243.25
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -8773.25
$ws.Range("N126").Value = -16938.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2729.1667
$ws.Range("J62").Value = 2413.0435
$ws.Range("L62").Value = 7239.130500000001
$ws.Range("N62").Value = -8611.130500000001
$ws.Range("H65").Value = 2729.1667
$ws.Range("J65").Value = 2413.0435
$ws.Range("L65").Value = 21717.3915
$ws.Range("N65").Value = -28581.3915
$ws.Range("H68").Value = 2678.9614
$ws.Range("J68").Value = 2781.7551
$ws.Range("L68").Value = 8345.265299999999
$ws.Range("N68").Value = -9967.265299999999
$ws.Range("H71").Value = 2678.9614
$ws.Range("J71").Value = 2781.7551
$ws.Range("L71").Value = 25035.7959
$ws.Range("N71").Value = -33147.7959
$ws.Range("H97").Value = 970.5
$ws.Range("I97").Value = 225
$ws.Range("K97").Value = 675
$ws.Range("M97").Value = -179
$ws.Range("H112").Value = 250001000
$ws.Range("J112").Value = 2000
$ws.Range("L112").Value = 6000
$ws.Range("N112").Value = -8216
$ws.Range("H115").Value = 1236.8
$ws.Range("J115").Value = 1515.5
$ws.Range("L115").Value = 4546.5
$ws.Range("N115").Value = -6896.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 400
$ws.Range("I14").Value = 400
$ws.Range("K14").Value = 400
$ws.Range("M14").Value = -232
$ws.Range("H80").Value = 5001.2856
$ws.Range("I80").Value = 5334.8335
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 5334.8335
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -4336.8335
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 5001.2856
$ws.Range("I83").Value = 5334.8335
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 26674.1675
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -21682.1675
$ws.Range("N83").Value = -24984
$ws.Range("H126").Value = 7488.3125
$ws.Range("I126").Value = 3642.8572
$ws.Range("K126").Value = 10928.5716
$ws.Range("M126").Value = -8458.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 16000
$ws.Range("J29").Value = 16000
$ws.Range("L29").Value = 16000
$ws.Range("N29").Value = -16590
$ws.Range("H40").Value = 12635.267
$ws.Range("J40").Value = 13200.8
$ws.Range("L40").Value = 13200.8
$ws.Range("N40").Value = -13472.8
$ws.Range("H55").Value = 2941904.8
$ws.Range("I55").Value = 3846417.5
$ws.Range("K55").Value = 3846417.5
$ws.Range("M55").Value = -3846244.5
$ws.Range("H61").Value = 4798.9
$ws.Range("I61").Value = 4727.857
$ws.Range("K61").Value = 4727.857
$ws.Range("M61").Value = -4525.857
$ws.Range("H82").Value = 6824.7334
$ws.Range("I82").Value = 6624.4287
$ws.Range("J82").Value = 7000
$ws.Range("K82").Value = 6624.4287
$ws.Range("L82").Value = 7000
$ws.Range("M82").Value = -6263.4287
$ws.Range("N82").Value = -7722
$ws.Range("H85").Value = 6824.7334
$ws.Range("I85").Value = 6624.4287
$ws.Range("J85").Value = 7000
$ws.Range("K85").Value = 6624.4287
$ws.Range("L85").Value = 7000
$ws.Range("M85").Value = -5376.4287
$ws.Range("N85").Value = -9496
$ws.Range("H100").Value = 4249.3076
$ws.Range("I100").Value = 2930.6365
$ws.Range("K100").Value = 2930.6365
$ws.Range("M100").Value = -2389.6365
$ws.Range("H113").Value = 4798.9
$ws.Range("I113").Value = 4727.857
$ws.Range("K113").Value = 4727.857
$ws.Range("M113").Value = -2557.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4902.2
$ws.Range("I132").Value = 2579.2307
$ws.Range("K132").Value = 7737.6921
$ws.Range("M132").Value = -5207.6921
$ws.Range("H136").Value = 3265.9375
$ws.Range("I136").Value = 1611.6111
$ws.Range("J136").Value = 5392.9287
$ws.Range("K136").Value = 4834.8333
$ws.Range("L136").Value = 16178.7861
$ws.Range("M136").Value = -2284.8333
$ws.Range("N136").Value = -21278.7861
